# The "Förändrad" (Changed) column C date was bumped for every data row
# (rows 2-484) from 2023-09-13 (45182) to 2023-09-15 (45184).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C484").Value = 45184
